$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 130
$ws.Range("F3").Value = 1321
$ws.Range("F4").Value = 1032
$ws.Range("F5").Value = 973
$ws.Range("F7").Value = 111
$ws.Range("F9").Value = 499
$ws.Range("F11").Value = 7
$ws.Range("F12").Value = 1869
$ws.Range("F13").Value = 4645
$ws.Range("F14").Value = 1335
$ws.Range("F15").Value = 132
$ws.Range("F16").Value = 2863
$ws.Range("F19").Value = 1163
$ws.Range("F20").Value = 3877
$ws.Range("F21").Value = 875
$ws.Range("F22").Value = 856
$ws.Range("F23").Value = 1562
$ws.Range("F24").Value = 54
$ws.Range("F25").Value = 2522
$ws.Range("F27").Value = 23
$ws.Range("F29").Value = 907
$ws.Range("F30").Value = 257
$ws.Range("F33").Value = 1004
$ws.Range("F34").Value = 273
$ws.Range("F35").Value = 58
$ws.Range("F37").Value = 107
$ws.Range("F38").Value = 1486
$ws.Range("F39").Value = 2042
$ws.Range("F40").Value = 970
$ws.Range("F41").Value = 22
$ws.Range("F42").Value = 31
$ws.Range("F43").Value = 539
$ws.Range("F44").Value = 149
$ws.Range("F45").Value = 624
$ws.Range("F46").Value = 333
$ws.Range("F47").Value = 168
$ws.Range("F48").Value = 180
$ws.Range("F49").Value = 98

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 26
$ws.Range("F12").Value = 133
$ws.Range("F20").Value = 25

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 583

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 583
$ws.Range("F3").Value = 1321
$ws.Range("F4").Value = 1032
$ws.Range("F5").Value = 973
$ws.Range("F7").Value = 111
$ws.Range("F8").Value = 26
$ws.Range("F9").Value = 26
$ws.Range("F12").Value = 499
$ws.Range("F13").Value = 7
$ws.Range("F14").Value = 4645
$ws.Range("F15").Value = 132
$ws.Range("F18").Value = 1163
$ws.Range("F19").Value = 3877
$ws.Range("F20").Value = 875
$ws.Range("F21").Value = 856
$ws.Range("F22").Value = 1562
$ws.Range("F23").Value = 54
$ws.Range("F24").Value = 2522
$ws.Range("F29").Value = 907
$ws.Range("F33").Value = 1004
$ws.Range("F34").Value = 273
$ws.Range("F36").Value = 1486
$ws.Range("F37").Value = 2042
$ws.Range("F39").Value = 970
$ws.Range("F40").Value = 22
$ws.Range("F42").Value = 31
$ws.Range("F44").Value = 539
$ws.Range("F45").Value = 624
$ws.Range("F46").Value = 333
$ws.Range("F47").Value = 168
$ws.Range("F48").Value = 180
$ws.Range("F49").Value = 98
$ws.Range("F50").Value = 25
